# Swap the deck's theme palette from "Integral" to the stock Office theme
# palette (ppt/theme/theme1.xml <a:clrScheme>). The font scheme ("Office" /
# Arial) and format scheme are already identical between the Integral theme
# and the Office theme in this deck, so only the 12 color-scheme entries
# need to change.
#
# PowerPoint's ColorScheme/.RGB property stores colors as a packed
# "0x00BBGGRR" long (Windows COLORREF order), so build each value from its
# R/G/B components with a small helper instead of hand-reversing hex
# strings.

function Convert-RGBtoBGRLong([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Target "Office Theme" color scheme, in ThemeColorScheme index order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5-10 Accent1-6, 11 Hyperlink, 12 FollowedHyperlink
$officeColors = @(
    @(0x00, 0x00, 0x00), # dk1
    @(0xFF, 0xFF, 0xFF), # lt1
    @(0x44, 0x54, 0x6A), # dk2
    @(0xE7, 0xE6, 0xE6), # lt2
    @(0x5B, 0x9B, 0xD5), # accent1
    @(0xED, 0x7D, 0x31), # accent2
    @(0xA5, 0xA5, 0xA5), # accent3
    @(0xFF, 0xC0, 0x00), # accent4
    @(0x44, 0x72, 0xC4), # accent5
    @(0x70, 0xAD, 0x47), # accent6
    @(0x05, 0x63, 0xC1), # hlink
    @(0x95, 0x4F, 0x72)  # folHlink
)

$p = $ppt.ActivePresentation

$slideTheme = $p.SlideMaster.Theme.ThemeColorScheme
$notesTheme = $p.NotesMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $rgb = $officeColors[$i - 1]
    $bgrLong = Convert-RGBtoBGRLong $rgb[0] $rgb[1] $rgb[2]
    $slideTheme.Item($i).RGB = $bgrLong
    $notesTheme.Item($i).RGB = $bgrLong
}
